$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the "Meta description" paragraph that used to sit right under the
#    document title (bold "Meta description" label + description text).
# ---------------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# ---------------------------------------------------------------------------
# 2. Insert a new bold paragraph right before the last paragraph of the
#    document ("Please create a feature image..." / soon to become the meta
#    description text), containing the page title text in bold.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertionPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)

$boldText = "Play Dragon Gate Trial for Free - Exciting Oriental-Themed Slot Game"
$xmlFragment = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>' + $boldText + '</w:t></w:r></w:p>'
$insertionPoint.InsertXML($xmlFragment)

# The inserted runs were merged onto the front of the old last paragraph, so
# split it back into two paragraphs right after the bold text we just added.
$mergedPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$splitOffset = $mergedPara.Range.Text.IndexOf("Please create a feature image")
$splitPoint = $d.Range($mergedPara.Range.Start + $splitOffset, $mergedPara.Range.Start + $splitOffset)
$splitPoint.InsertParagraphBefore()

# ---------------------------------------------------------------------------
# 3. Replace the old "Please create a feature image..." instructions with the
#    meta description text (keeping the paragraph's italic formatting).
# ---------------------------------------------------------------------------
$oldText = "Please create a feature image for Dragon Gate Trial that fits the following criteria: - The image should be in cartoon style - The image should feature a happy Maya warrior with glasses. Note: The image should not include any references to the game's actual theme of Chinese New Year and dragons. The Maya warrior should be the main focus of the image."
$newText = "Discover the winning possibilities and enticing combination of Asian themes and classic slots when you play Dragon Gate Trial for free."

$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $newText, 2) | Out-Null
